$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 1849.5
$ws.Cells.Item(12, 9).Value = 1798.6666
$ws.Cells.Item(12, 10).Value = 2002
$ws.Cells.Item(12, 11).Value = 1798.6666
$ws.Cells.Item(12, 12).Value = 2002
$ws.Cells.Item(12, 13).Value = -1628.6666
$ws.Cells.Item(12, 14).Value = -2342
$ws.Cells.Item(51, 8).Value = 12499.25
$ws.Cells.Item(51, 9).Value = 9999
$ws.Cells.Item(51, 11).Value = 9999
$ws.Cells.Item(51, 13).Value = -9515
$ws.Cells.Item(55, 8).Value = 134.88235
$ws.Cells.Item(55, 9).Value = 106.77778
$ws.Cells.Item(55, 10).Value = 166.5
$ws.Cells.Item(55, 11).Value = 106.77778
$ws.Cells.Item(55, 12).Value = 166.5
$ws.Cells.Item(55, 13).Value = 107.22222
$ws.Cells.Item(55, 14).Value = -594.5
$ws.Cells.Item(64, 8).Value = 3880
$ws.Cells.Item(64, 9).Value = 3880
$ws.Cells.Item(64, 11).Value = 3880
$ws.Cells.Item(64, 13).Value = -3632
$ws.Cells.Item(67, 8).Value = 3880
$ws.Cells.Item(67, 9).Value = 3880
$ws.Cells.Item(67, 11).Value = 3880
$ws.Cells.Item(67, 13).Value = -3022
$ws.Cells.Item(70, 8).Value = 2010.421
$ws.Cells.Item(70, 9).Value = 1599
$ws.Cells.Item(70, 10).Value = 2058.8235
$ws.Cells.Item(70, 11).Value = 4797
$ws.Cells.Item(70, 12).Value = 6176.470499999999
$ws.Cells.Item(70, 13).Value = -4527
$ws.Cells.Item(70, 14).Value = -6716.470499999999
$ws.Cells.Item(73, 8).Value = 2010.421
$ws.Cells.Item(73, 9).Value = 1599
$ws.Cells.Item(73, 10).Value = 2058.8235
$ws.Cells.Item(73, 11).Value = 4797
$ws.Cells.Item(73, 12).Value = 6176.470499999999
$ws.Cells.Item(73, 13).Value = -3861
$ws.Cells.Item(73, 14).Value = -8048.470499999999
$ws.Cells.Item(112, 8).Value = 3227.923
$ws.Cells.Item(112, 10).Value = 3227.923
$ws.Cells.Item(112, 12).Value = 9683.769
$ws.Cells.Item(112, 14).Value = -11899.769
$ws.Cells.Item(116, 8).Value = 4095.0715
$ws.Cells.Item(116, 10).Value = 4310.1665
$ws.Cells.Item(116, 12).Value = 4310.1665
$ws.Cells.Item(116, 14).Value = -11194.1665
$ws.Cells.Item(132, 8).Value = 2115.0833
$ws.Cells.Item(132, 9).Value = 2034.6364
$ws.Cells.Item(132, 11).Value = 6103.9092
$ws.Cells.Item(132, 13).Value = -3573.9092

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3983.457
$ws.Cells.Item(32, 9).Value = 4476.067
$ws.Cells.Item(32, 11).Value = 4476.067
$ws.Cells.Item(32, 13).Value = -4189.067
$ws.Cells.Item(61, 8).Value = 3490.5557
$ws.Cells.Item(61, 9).Value = 2401.4167
$ws.Cells.Item(61, 10).Value = 5668.8335
$ws.Cells.Item(61, 11).Value = 2401.4167
$ws.Cells.Item(61, 12).Value = 5668.8335
$ws.Cells.Item(61, 13).Value = -2189.4167
$ws.Cells.Item(61, 14).Value = -6092.8335
$ws.Cells.Item(74, 8).Value = 1751.5555
$ws.Cells.Item(74, 9).Value = 1720.5
$ws.Cells.Item(74, 11).Value = 1720.5
$ws.Cells.Item(74, 13).Value = -846.5
$ws.Cells.Item(77, 8).Value = 1751.5555
$ws.Cells.Item(77, 9).Value = 1720.5
$ws.Cells.Item(77, 11).Value = 8602.5
$ws.Cells.Item(77, 13).Value = -4234.5
$ws.Cells.Item(136, 8).Value = 3490.5557
$ws.Cells.Item(136, 9).Value = 2401.4167
$ws.Cells.Item(136, 10).Value = 5668.8335
$ws.Cells.Item(136, 11).Value = 7204.250100000001
$ws.Cells.Item(136, 12).Value = 17006.5005
$ws.Cells.Item(136, 13).Value = -4654.250100000001
$ws.Cells.Item(136, 14).Value = -22106.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(24, 8).Value = 6304
$ws.Cells.Item(24, 9).Value = 6304
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 11).Value = 6304
$ws.Cells.Item(24, 12).Value = 0
$ws.Cells.Item(24, 13).ClearContents()
$ws.Cells.Item(24, 14).Value = -6069
$ws.Cells.Item(36, 8).Value = 3041
$ws.Cells.Item(36, 9).Value = 0
$ws.Cells.Item(36, 11).Value = 0
$ws.Cells.Item(36, 13).ClearContents()
$ws.Cells.Item(134, 8).Value = 1098
$ws.Cells.Item(134, 9).Value = 1098
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 3294
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 13).ClearContents()
$ws.Cells.Item(134, 14).Value = -759

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 5849.6665
$ws.Cells.Item(22, 9).Value = 799
$ws.Cells.Item(22, 10).Value = 8375
$ws.Cells.Item(22, 11).Value = 799
$ws.Cells.Item(22, 12).Value = 8375
$ws.Cells.Item(22, 13).Value = -449
$ws.Cells.Item(22, 14).Value = -9075
$ws.Cells.Item(31, 8).Value = 1199
$ws.Cells.Item(31, 9).Value = 1199
$ws.Cells.Item(31, 11).Value = 1199
$ws.Cells.Item(31, 13).Value = -904
$ws.Cells.Item(34, 8).Value = 1199
$ws.Cells.Item(34, 9).Value = 1199
$ws.Cells.Item(34, 11).Value = 1199
$ws.Cells.Item(34, 13).Value = -997

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 533
$ws.Cells.Item(92, 9).Value = 382.83334
$ws.Cells.Item(92, 11).Value = 1148.50002
$ws.Cells.Item(92, 13).Value = 99.49998000000005
$ws.Cells.Item(106, 8).Value = 0
$ws.Cells.Item(106, 10).Value = 0
$ws.Cells.Item(106, 12).ClearContents()
$ws.Cells.Item(106, 14).Value = 0
$ws.Cells.Item(131, 8).Value = 756.1667
$ws.Cells.Item(134, 8).Value = 1073.8
$ws.Cells.Item(134, 9).Value = 1073.8
$ws.Cells.Item(134, 11).Value = 3221.4
$ws.Cells.Item(134, 13).Value = 1848.6
$ws.Cells.Item(139, 8).Value = 4805.8335
$ws.Cells.Item(139, 9).Value = 3767
$ws.Cells.Item(139, 13).Value = -6161

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(93, 8).Value = 50000
$ws.Cells.Item(93, 10).Value = 50000
$ws.Cells.Item(93, 12).Value = 50000
$ws.Cells.Item(93, 14).Value = -53744
$ws.Cells.Item(126, 8).Value = 8168.5
$ws.Cells.Item(126, 9).Value = 6002.75
$ws.Cells.Item(126, 11).Value = 18008.25
$ws.Cells.Item(126, 13).Value = -15538.25
$ws.Cells.Item(132, 8).Value = 2111.4
$ws.Cells.Item(132, 9).Value = 1845.6923
$ws.Cells.Item(132, 10).Value = 2604.8572
$ws.Cells.Item(132, 11).Value = 5537.0769
$ws.Cells.Item(132, 12).Value = 7814.571599999999
$ws.Cells.Item(132, 13).Value = -3007.0769
$ws.Cells.Item(132, 14).Value = -12874.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1000
$ws.Cells.Item(16, 9).Value = 1000
$ws.Cells.Item(16, 11).Value = 1000
$ws.Cells.Item(16, 13).Value = -830
$ws.Cells.Item(22, 8).Value = 2054
$ws.Cells.Item(22, 9).Value = 1578.8
$ws.Cells.Item(22, 10).Value = 2450
$ws.Cells.Item(22, 11).Value = 1578.8
$ws.Cells.Item(22, 12).Value = 2450
$ws.Cells.Item(22, 13).Value = -1283.8
$ws.Cells.Item(22, 14).Value = -3040
$ws.Cells.Item(27, 8).Value = 2054
$ws.Cells.Item(27, 9).Value = 1578.8
$ws.Cells.Item(27, 10).Value = 2450
$ws.Cells.Item(27, 11).Value = 1578.8
$ws.Cells.Item(27, 12).Value = 2450
$ws.Cells.Item(27, 13).Value = -1471.8
$ws.Cells.Item(27, 14).Value = -2664
$ws.Cells.Item(46, 8).Value = 4172.773
$ws.Cells.Item(46, 10).Value = 4800
$ws.Cells.Item(46, 12).Value = 4800
$ws.Cells.Item(46, 14).Value = -5176
$ws.Cells.Item(95, 8).Value = 45167
$ws.Cells.Item(95, 10).Value = 45167
$ws.Cells.Item(95, 12).Value = 45167
$ws.Cells.Item(95, 14).Value = -50659
$ws.Cells.Item(132, 8).Value = 2362.6365
$ws.Cells.Item(132, 9).Value = 2373.75
$ws.Cells.Item(132, 10).Value = 2333
$ws.Cells.Item(132, 11).Value = 7121.25
$ws.Cells.Item(132, 12).Value = 6999
$ws.Cells.Item(132, 13).Value = -4591.25
$ws.Cells.Item(132, 14).Value = -12059

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 1171.6666
$ws.Cells.Item(96, 9).Value = 1008.4
$ws.Cells.Item(96, 11).Value = 1008.4
$ws.Cells.Item(96, 13).Value = 364.6
$ws.Cells.Item(136, 8).Value = 845.5
$ws.Cells.Item(136, 9).Value = 500.625
$ws.Cells.Item(136, 11).Value = 1501.875
$ws.Cells.Item(136, 13).Value = 1048.125
